$d = $word.ActiveDocument

$d.Content.Find.Execute("204×3=612", $true, $false, $false, $false, $false, $true, 1, $false, "255×4=1020", 2)
$d.Content.Find.Execute("509×5=2545", $true, $false, $false, $false, $false, $true, 1, $false, "267×4=1068", 2)
$d.Content.Find.Execute("268×4=1072", $true, $false, $false, $false, $false, $true, 1, $false, "217×6=1302", 2)
$d.Content.Find.Execute("537×4=2148", $true, $false, $false, $false, $false, $true, 1, $false, "954×9=8586", 2)
$d.Content.Find.Execute("993×2=1986", $true, $false, $false, $false, $false, $true, 1, $false, "529×6=3174", 2)
$d.Content.Find.Execute("492×8=3936", $true, $false, $false, $false, $false, $true, 1, $false, "261×5=1305", 2)
$d.Content.Find.Execute("767×5=3835", $true, $false, $false, $false, $false, $true, 1, $false, "672×3=2016", 2)
$d.Content.Find.Execute("424×4=1696", $true, $false, $false, $false, $false, $true, 1, $false, "113×3=339", 2)
$d.Content.Find.Execute("797×6=4782", $true, $false, $false, $false, $false, $true, 1, $false, "318×8=2544", 2)
$d.Content.Find.Execute("441×5=2205", $true, $false, $false, $false, $false, $true, 1, $false, "222×6=1332", 2)
$d.Content.Find.Execute("607×2=1214", $true, $false, $false, $false, $false, $true, 1, $false, "345×4=1380", 2)
$d.Content.Find.Execute("265×3=795", $true, $false, $false, $false, $false, $true, 1, $false, "101×5=505", 2)
$d.Content.Find.Execute("720×4=2880", $true, $false, $false, $false, $false, $true, 1, $false, "198×5=990", 2)
$d.Content.Find.Execute("473×8=3784", $true, $false, $false, $false, $false, $true, 1, $false, "846×7=5922", 2)
$d.Content.Find.Execute("299×4=1196", $true, $false, $false, $false, $false, $true, 1, $false, "723×8=5784", 2)
$d.Content.Find.Execute("843×4=3372", $true, $false, $false, $false, $false, $true, 1, $false, "580×6=3480", 2)
$d.Content.Find.Execute("148×3=444", $true, $false, $false, $false, $false, $true, 1, $false, "484×3=1452", 2)
$d.Content.Find.Execute("194×7=1358", $true, $false, $false, $false, $false, $true, 1, $false, "394×8=3152", 2)
$d.Content.Find.Execute("393×8=3144", $true, $false, $false, $false, $false, $true, 1, $false, "177×9=1593", 2)
$d.Content.Find.Execute("492×3=1476", $true, $false, $false, $false, $false, $true, 1, $false, "609×5=3045", 2)
$d.Content.Find.Execute("930×7=6510", $true, $false, $false, $false, $false, $true, 1, $false, "202×4=808", 2)
$d.Content.Find.Execute("692×9=6228", $true, $false, $false, $false, $false, $true, 1, $false, "206×5=1030", 2)
$d.Content.Find.Execute("348×6=2088", $true, $false, $false, $false, $false, $true, 1, $false, "212×6=1272", 2)
$d.Content.Find.Execute("388×9=3492", $true, $false, $false, $false, $false, $true, 1, $false, "231×7=1617", 2)
$d.Content.Find.Execute("668×4=2672", $true, $false, $false, $false, $false, $true, 1, $false, "204×2=408", 2)
